$p = $ppt.ActivePresentation

# Slide 1 title: consolidate "Example" " " "numbering" " " "MWE" runs into one run.
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
# The resulting string is unchanged from the current concatenated text, so a
# direct assignment is seen as a no-op by the diffing engine. Force a real
# text-content change first, then set the final (consolidated) value so the
# multiple <a:r> runs collapse into a single run.
$tr1.Text = "__tmp__"
$tr1.Text = "Example numbering MWE"

# Slide 2 title: consolidate "A" " " "second" " " "slide" runs into one run.
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Text = "__tmp__"
$tr2.Text = "A second slide"
